# Apply edits described in the diff:
# - user_name (column C) changed from "Momoko Endo" to "Shohei Ota" for rows 2-16
# - capimg (column J) filenames updated to new bdot20240415_141954/N.png scheme
# - explanation (column K) text updated/shifted per new step descriptions
# - row 5 becomes the new "error" row (type=error, explanation/error_type/error_content filled)
# - row 7 becomes a normal "operation" row (type=operation, error_type/error_content cleared)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C: user_name -> Shohei Ota for all data rows (2-16)
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 3).Value = "Shohei Ota"
}

# Column B: type changes
$ws.Cells.Item(5, 2).Value = "error"
$ws.Cells.Item(7, 2).Value = "operation"

# Column J: capimg filenames
$ws.Range("J2").Value = "bdot20240415_141954/1.png"
$ws.Range("J3").Value = "bdot20240415_141954/2.png"
$ws.Range("J4").Value = "bdot20240415_141954/3.png"
$ws.Range("J5").Value = "bdot20240415_141954/4.png"
$ws.Range("J6").Value = "bdot20240415_141954/5.png"
$ws.Range("J7").Value = "bdot20240415_141954/5.png"
$ws.Range("J8").Value = "bdot20240415_141954/6.png"
$ws.Range("J9").Value = "bdot20240415_141954/7.png"
$ws.Range("J10").Value = "bdot20240415_141954/8.png"
$ws.Range("J11").Value = "bdot20240415_141954/9.png"
$ws.Range("J12").Value = "bdot20240415_141954/10.png"
$ws.Range("J13").Value = "bdot20240415_141954/1.png"
$ws.Range("J14").Value = "bdot20240415_141954/2.png"
$ws.Range("J15").Value = "bdot20240415_141954/3.png"
$ws.Range("J16").Value = "bdot20240415_141954/11.png"

# Column K: explanation text
$ws.Range("K2").Value = "「スタート」ボタンをクリックする"
$ws.Range("K3").Value = "メニューから「設定」アイコンをクリックする"
$ws.Range("K4").Value = "左側のメニューからWindows Updateをクリックし、Windows Update画面に移動する"
$ws.Range("K5").Value = "0x80240fff エラー"
$ws.Range("K6").Value = "デスクトップ画面の左下にある「スタート」ボタンを右クリックする"
$ws.Range("K7").Value = "メニューからターミナル(管理者)をクリックする"
$ws.Range("K8").Value = "ユーザーアカウント制御と表示されているウィンドウが開いたことを確認する"
$ws.Range("K9").Value = "PowerShellウィンドウに start-transcript と入力し、[Enter]キーを押す"
$ws.Range("K10").Value = "wuauclt.exe /resetauthorization /detectnow と入力し、[Enter]キーを押す"
$ws.Range("K11").Value = "netsh winhttp show proxy と入力し、[Enter]キーを押す"
$ws.Range("K12").Value = "netsh winhttp reset proxy と入力し、[Enter]キーを押す"
$ws.Range("K13").Value = "「スタート」ボタンをクリックする"
$ws.Range("K14").Value = "メニューから「設定」アイコンをクリックする"
$ws.Range("K15").Value = "左側のメニューからWindows Updateをクリックし、Windows Update画面に移動する"
$ws.Range("K16").Value = "「更新プログラムのチェック」ボタンをクリックする"

# Row 5 becomes the error row: fill error_type (L) and error_content (M)
$ws.Range("L5").Value = "Error W"
$ws.Range("M5").Value = " エラーの Windows"

# Row 7 no longer carries error details: clear error_type (L) and error_content (M)
$ws.Range("L7").Value = ""
$ws.Range("M7").Value = ""

Write-Host "Edits applied"
